# Auto-generated edit script: update cryptos list values per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.536.23"
$ws.Range("E2").Value = "  +0.19%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.107.79"
$ws.Range("E3").Value = "  +0.63%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.012"
$ws.Range("E4").Value = "  +0.77%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "335.81"
$ws.Range("E5").Value = "  +1.80%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.010"
$ws.Range("E6").Value = "  +0.71%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5240"
$ws.Range("E7").Value = "  +0.52%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4530"
$ws.Range("E8").Value = "  +3.88%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "55.54"
$ws.Range("E9").Value = "  +1.85%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09019"
$ws.Range("E10").Value = "  +1.67%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.168"
$ws.Range("E11").Value = "  +1.16%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.55"
$ws.Range("E12").Value = "  +0.65%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.115.19"
$ws.Range("E13").Value = "  +0.98%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.831"
$ws.Range("E14").Value = "  +2.38%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.108"
$ws.Range("E15").Value = "  +5.73%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001171"
$ws.Range("E16").Value = "  +4.37%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "96.97"
$ws.Range("E17").Value = "  +1.18%  "
$ws.Range("E18").Value = "  +0.92%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06688"
$ws.Range("E19").Value = "  +1.52%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.32"
$ws.Range("E20").Value = "  +0.15%  "
$ws.Range("E21").Value = "  +0.76%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.230"
$ws.Range("E22").Value = "  -0.29%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.597.97"
$ws.Range("E23").Value = "  +0.25%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.80"
$ws.Range("E24").Value = "  +4.47%  "
$ws.Range("E25").Value = "  +0.72%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.357.97"
$ws.Range("E26").Value = "  +0.79%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.34"
$ws.Range("E27").Value = "  +0.39%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "163.28"
$ws.Range("E28").Value = "  +0.43%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.504"
$ws.Range("E29").Value = "  -2.12%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.39"
$ws.Range("E30").Value = "  +1.24%  "
$ws.Range("E31").Value = "  +2.44%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1066"
$ws.Range("E32").Value = "  -0.08%  "
$ws.Range("B33").Value = "ARBITRUM"
$ws.Range("C33").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.626"
$ws.Range("E33").Value = "  -0.72%  "
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.333"
$ws.Range("E34").Value = "  +3.17%  "
$ws.Range("E35").Value = "  +1.81%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.39"
$ws.Range("E36").Value = "  +3.87%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.869"
$ws.Range("E37").Value = "  +7.88%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02609"
$ws.Range("E38").Value = "  +1.42%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06798"
$ws.Range("E39").Value = "  -0.32%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2307"
$ws.Range("E40").Value = "  +2.36%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.56"
$ws.Range("E41").Value = "  -1.41%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6821"
$ws.Range("E42").Value = "  -0.81%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.257"
$ws.Range("E43").Value = "  -0.04%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6421"
$ws.Range("E44").Value = "  +0.95%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.00"
$ws.Range("E45").Value = "  +1.03%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.300"
$ws.Range("E46").Value = "  +4.76%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.668"
$ws.Range("E47").Value = "  +1.18%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.00000000355"
$ws.Range("E48").Value = "  +19.10%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.248"
$ws.Range("E49").Value = "  +0.64%  "
$ws.Range("E50").Value = "  -2.18%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "83.12"
$ws.Range("E51").Value = "  +1.65%  "
